$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix "Panner" -> "Paneer" typo in the two menu item names
$ws.Range("A12").Value = "7 Inch Pizza Paneer"
$ws.Range("A8").Value = "Pizza Paneer, Veggie ( Onion and Capsicum and corn) 10'"

# Price updates
$ws.Range("C5").Value = 120
$ws.Range("C6").Value = 140
$ws.Range("C7").Value = 150

# Update view: scroll back to top-left and move selection to A8
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A8").Select()

